$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "89.616.18"
$ws.Range("E2").Value = "  -2.06%  "

$ws.Range("D3").Value = "3.100.06"
$ws.Range("E3").Value = "  -2.78%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.26"
$ws.Range("E5").Value = "  -1.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "623.06"
$ws.Range("E6").Value = "  -1.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.374"
$ws.Range("E7").Value = "  -6.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.827"
$ws.Range("E8").Value = "  +15.86%  "

$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").Value = "3.098.23"
$ws.Range("E10").Value = "  -2.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.610"
$ws.Range("E11").Value = "  +7.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.183"
$ws.Range("E12").Value = "  +1.50%  "

$ws.Range("E13").Value = "  -4.89%  "

$ws.Range("E14").Value = "  +0.02%  "

$ws.Range("D15").Value = "89.277.75"
$ws.Range("E15").Value = "  -1.91%  "

$ws.Range("D16").Value = "3.673.62"
$ws.Range("E16").Value = "  -2.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "32.33"
$ws.Range("E17").Value = "  -0.35%  "

$ws.Range("D18").Value = "3.105.29"
$ws.Range("E18").Value = "  -2.13%  "

$ws.Range("E19").Value = "  +2.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000216"
$ws.Range("E20").Value = "  +1.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.48"
$ws.Range("E21").Value = "  +1.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "426.61"
$ws.Range("E22").Value = "  -2.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.29"
$ws.Range("E23").Value = "  -2.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.96"
$ws.Range("E24").Value = "  -0.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.52"
$ws.Range("E25").Value = "  +5.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.07"
$ws.Range("E26").Value = "  +2.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "84.14"
$ws.Range("E27").Value = "  +4.66%  "

$ws.Range("D28").Value = "3.259.96"
$ws.Range("E28").Value = "  -2.26%  "

$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.09"
$ws.Range("E30").Value = "  +8.86%  "

$ws.Range("E31").Value = "  +5.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.15"
$ws.Range("E32").Value = "  -1.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "511.00"
$ws.Range("E33").Value = "  -1.62%  "

$ws.Range("E34").Value = "  -8.00%  "

$ws.Range("E35").Value = "  -3.64%  "

$ws.Range("E36").Value = "  -1.87%  "

$ws.Range("E37").Value = "  -3.97%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.46"
$ws.Range("E38").Value = "  +0.51%  "

$ws.Range("E39").Value = "  -0.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.130"
$ws.Range("E40").Value = "  +2.70%  "

$ws.Range("E41").Value = "  +0.08%  "

$ws.Range("E43").Value = "  -1.28%  "

$ws.Range("E44").Value = "  -4.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.135"
$ws.Range("E45").Value = "  +7.52%  "

$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0719"
$ws.Range("E46").Value = "  +17.54%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "145.73"
$ws.Range("E47").Value = "  -0.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "43.48"
$ws.Range("E48").Value = "  -0.75%  "

$ws.Range("E49").Value = "  +1.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "160.69"
$ws.Range("E50").Value = "  -6.00%  "

$ws.Range("E51").Value = "  -4.77%  "
